# Commit: "change country USA to United States"
#
# The port/country lookup table lists "USA" as the country for 14 ports.
# This edit renames "USA" to "United States" for all of those ports except
# the "Reserve" entry (row 120), which keeps its original "USA" value.
#
# Affected rows (port -> country cell in column B):
#   20  Baton Rouge
#   22  Bayport
#   23  Beaumont
#   41  Corpus Christi
#   61  Freeport
#   63  Galveston
#   68  Houston
#   78  LOOP Terminal
#  101  Norco
#  109  Port Arthur
#  113  Port Sulphur
#  138  St James
#  149  Texas City
#
# Row 120 (Reserve) is intentionally left untouched -> stays "USA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToRename = @(20, 22, 23, 41, 61, 63, 68, 78, 101, 109, 113, 138, 149)

foreach ($r in $rowsToRename) {
    $ws.Cells.Item($r, 2).Value = "United States"
}

# Keep gridlines visible (workbook default) through the save round-trip.
$excel.ActiveWindow.DisplayGridlines = $true

# Match the workbook's final view/selection state recorded in the edit
# (scrolled back to the top, with B150 as the active cell).
[void]$ws.Range("A1").Select()
[void]$ws.Range("B150").Select()
